# #3473 replaced two properties that had gaps
#
# Two properties (row 2: "Medstar POB North Tower" -> "Medstar POB South
# Tower", and row 10: "DPW Vehicle Maintenance Facility 2" -> "School
# Without Walls @ Francis Stevens") were fully replaced with new property
# records, and several other rows had small data corrections (addresses,
# owner names, postal codes, gross area). The "Year Built" column (I) also
# had its accidental date number-format removed so it shows plain years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Medstar POB North Tower -> Medstar POB South Tower ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319
$ws.Range("M2").Value = 11.1
$ws.Range("N2").Value = 121.2

# --- Row 3: 1801 Pennsylvania Ave. -> 1801 Pennsylvania Avenue, LLC ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

# --- Row 4: GSA: 300 E Street SW corrections ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# --- Row 5: Paul H.Nitze gross area correction ---
$ws.Range("L5").Value = 58717

# --- Row 6: President Madison Apartments -> Hampton House ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580
$ws.Range("M6").Value = 3.4
$ws.Range("N6").Value = 58.9

# --- Row 7: 3303 Water Street corrections ---
$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

# --- Row 8: 15th and H Street Associates LLP address correction ---
$ws.Range("E8").Value = "1428 H ST NW"

# --- Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991
$ws.Range("M10").Value = 4.2
$ws.Range("N10").Value = 70.5
$ws.Range("P10").Value = 68

# The "Year Built" column had a stray date number format (it was showing
# years like 1967 as dates); clear that formatting so the raw year values
# display as plain numbers, matching every other data column.
$ws.Range("I2:I10").ClearFormats()

# Restore the selection to what was left after the edits.
$ws.Range("N14").Select()
